# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.211.11"
$ws.Range("E2").Value = "  +13.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.675.43"
$ws.Range("E3").Value = "  +8.24%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.38"
$ws.Range("E5").Value = "  +9.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3737"
$ws.Range("E7").Value = "  +2.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3451"
$ws.Range("E8").Value = "  +7.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.60"
$ws.Range("E9").Value = "  +15.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.183"
$ws.Range("E10").Value = "  +6.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07310"
$ws.Range("E11").Value = "  +5.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.0000"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("E13").Value = "  +7.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.115"
$ws.Range("E14").Value = "  +6.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.776"
$ws.Range("E15").Value = "  +5.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.676.02"
$ws.Range("E16").Value = "  +8.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001111"
$ws.Range("E17").Value = "  +5.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9993"
$ws.Range("E18").Value = "  +3.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06733"
$ws.Range("E19").Value = "  +9.32%  "
$ws.Range("E20").Value = "  +11.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.57"
$ws.Range("E21").Value = "  +7.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.147"
$ws.Range("E22").Value = "  +6.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.04"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.177.50"
$ws.Range("E24").Value = "  +13.09%  "
$ws.Range("E25").Value = "  +3.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.364"
$ws.Range("E26").Value = "  -9.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.665"
$ws.Range("E27").Value = "  +16.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.58"
$ws.Range("E28").Value = "  +2.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.55"
$ws.Range("E29").Value = "  +9.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.862.22"
$ws.Range("E30").Value = "  +8.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.79"
$ws.Range("E31").Value = "  +6.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.501"
$ws.Range("E32").Value = "  +23.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.096"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9932"
$ws.Range("E34").Value = "  +11.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.770"
$ws.Range("E35").Value = "  +14.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08519"
$ws.Range("E36").Value = "  +5.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.65"
$ws.Range("E37").Value = "  +15.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06508"
$ws.Range("E38").Value = "  +10.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.378"
$ws.Range("E39").Value = "  +7.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.927"
$ws.Range("E40").Value = "  +12.66%  "
$ws.Range("E41").Value = "  +10.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.278"
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2145"
$ws.Range("E43").Value = "  +6.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6205"
$ws.Range("E44").Value = "  +12.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9987"
$ws.Range("E45").Value = "  +2.99%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.811"
$ws.Range("E46").Value = "  +6.29%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.17"
$ws.Range("E47").Value = "  +4.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5955"
$ws.Range("E48").Value = "  +7.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.48"
$ws.Range("E49").Value = "  +4.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.034"
$ws.Range("E50").Value = "  +7.80%  "
$ws.Range("E51").Value = "  +8.15%  "
